# Rename the inline logo pictures living in the document's headers and
# footers. The BTEC logo (header, both "default" and "first page" headers)
# swaps from image1.jpg -> image2.jpg, and the Pearson/Edexcel logo
# (footer, both "default" and "first page" footers) swaps from
# image2.png -> image1.png. Alt-text / picture content is untouched -
# only the shape's OOXML `name` changes.

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    # Headers: BTec_Logo-Orange picture, image1.jpg -> image2.jpg
    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            for ($i = 1; $i -le $header.Range.InlineShapes.Count; $i++) {
                $shape = $header.Range.InlineShapes.Item($i)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape.Name = "image2.jpg"
                }
            }
        }
    }

    # Footers: Pearson Edexcel logo picture, image2.png -> image1.png
    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            for ($i = 1; $i -le $footer.Range.InlineShapes.Count; $i++) {
                $shape = $footer.Range.InlineShapes.Item($i)
                if ($shape.AlternativeText -like "*PearsonLogo.png") {
                    $shape.Name = "image1.png"
                }
            }
        }
    }
}
